$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 1518.0714  # H98: 1503.3103 -> 1518.0714
$ws.Cells.Item(98, 9).Value = 1444.6666  # I98: 1432 -> 1444.6666
$ws.Cells.Item(98, 11).Value = 1444.6666  # K98: 1432 -> 1444.6666
$ws.Cells.Item(98, 13).Value = 53.33339999999998  # M98: 66 -> 53.33339999999998
$ws.Cells.Item(122, 8).Value = 1518.0714  # H122: 1503.3103 -> 1518.0714
$ws.Cells.Item(122, 9).Value = 1444.6666  # I122: 1432 -> 1444.6666
$ws.Cells.Item(122, 11).Value = 4333.9998  # K122: 4296 -> 4333.9998
$ws.Cells.Item(122, 13).Value = -1883.9998  # M122: -1846 -> -1883.9998
$ws.Cells.Item(132, 8).Value = 4741.1  # H132: 4894.483 -> 4741.1
$ws.Cells.Item(132, 9).Value = 4741.1  # I132: 4894.483 -> 4741.1
$ws.Cells.Item(132, 11).Value = 14223.3  # K132: 14683.449 -> 14223.3
$ws.Cells.Item(132, 13).Value = -11693.3  # M132: -12153.449 -> -11693.3
$ws.Cells.Item(137, 8).Value = 5572849  # H137: 4180302.5 -> 5572849
$ws.Cells.Item(137, 9).Value = 10001986  # I137: 6252240 -> 10001986
$ws.Cells.Item(137, 11).Value = 30005958  # K137: 18756720 -> 30005958
$ws.Cells.Item(137, 13).Value = -30003408  # M137: -18754170 -> -30003408
$ws.Cells.Item(138, 8).Value = 6609.12  # H138: 6436.675 -> 6609.12
$ws.Cells.Item(138, 10).Value = 7379.7803  # J138: 7405.871 -> 7379.7803
$ws.Cells.Item(138, 12).Value = 22139.3409  # L138: 22217.613 -> 22139.3409
$ws.Cells.Item(138, 14).Value = -32419.3409  # N138: -32497.613 -> -32419.3409

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1541116.2  # H32: 1304029.4 -> 1541116.2
$ws.Cells.Item(32, 9).Value = 773632.5  # I32: 651485.6 -> 773632.5
$ws.Cells.Item(32, 10).Value = 6803862  # J32: 5953404 -> 6803862
$ws.Cells.Item(32, 11).Value = 773632.5  # K32: 651485.6 -> 773632.5
$ws.Cells.Item(32, 12).Value = 6803862  # L32: 5953404 -> 6803862
$ws.Cells.Item(32, 13).Value = -773345.5  # M32: -651198.6 -> -773345.5
$ws.Cells.Item(32, 14).Value = -6804436  # N32: -5953978 -> -6804436
$ws.Cells.Item(56, 8).Value = 0  # H56: 15000 -> 0
$ws.Cells.Item(56, 9).Value = 0  # I56: 15000 -> 0
$ws.Cells.Item(56, 11).Value = 0  # K56: 15000 -> 0
$ws.Cells.Item(56, 13).ClearContents()  # M56: -14258 -> (removed)
$ws.Cells.Item(63, 8).Value = 2659.2  # H63: 2220.2856 -> 2659.2
$ws.Cells.Item(63, 9).Value = 2324  # I63: 2049.2 -> 2324
$ws.Cells.Item(63, 10).Value = 4000  # J63: 2648 -> 4000
$ws.Cells.Item(63, 11).Value = 2324  # K63: 2049.2 -> 2324
$ws.Cells.Item(63, 12).Value = 4000  # L63: 2648 -> 4000
$ws.Cells.Item(63, 13).Value = -1638  # M63: -1363.2 -> -1638
$ws.Cells.Item(63, 14).Value = -5372  # N63: -4020 -> -5372
$ws.Cells.Item(66, 8).Value = 2659.2  # H66: 2220.2856 -> 2659.2
$ws.Cells.Item(66, 9).Value = 2324  # I66: 2049.2 -> 2324
$ws.Cells.Item(66, 10).Value = 4000  # J66: 2648 -> 4000
$ws.Cells.Item(66, 11).Value = 11620  # K66: 10246 -> 11620
$ws.Cells.Item(66, 12).Value = 20000  # L66: 13240 -> 20000
$ws.Cells.Item(66, 13).Value = -8188  # M66: -6814 -> -8188
$ws.Cells.Item(66, 14).Value = -26864  # N66: -20104 -> -26864
$ws.Cells.Item(132, 8).Value = 1755.8572  # H132: 1715.9828 -> 1755.8572
$ws.Cells.Item(132, 9).Value = 1540.234  # I132: 1501.8368 -> 1540.234
$ws.Cells.Item(132, 11).Value = 4620.701999999999  # K132: 4505.5104 -> 4620.701999999999
$ws.Cells.Item(132, 13).Value = -2090.701999999999  # M132: -1975.5104 -> -2090.701999999999

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(35, 8).Value = 43291  # H35: 44999.5 -> 43291
$ws.Cells.Item(35, 10).Value = 43291  # J35: 44999.5 -> 43291
$ws.Cells.Item(35, 12).Value = 43291  # L35: 44999.5 -> 43291
$ws.Cells.Item(35, 14).Value = -43911  # N35: -45619.5 -> -43911
$ws.Cells.Item(94, 8).Value = 95239860  # H94: 102565920 -> 95239860
$ws.Cells.Item(94, 9).Value = 95239860  # I94: 102565920 -> 95239860
$ws.Cells.Item(94, 11).Value = 95239860  # K94: 102565920 -> 95239860
$ws.Cells.Item(94, 13).Value = -95239409  # M94: -102565469 -> -95239409
$ws.Cells.Item(132, 8).Value = 99636.37  # H132: 99076.84 -> 99636.37
$ws.Cells.Item(132, 10).Value = 99636.37  # J132: 99076.84 -> 99636.37
$ws.Cells.Item(132, 12).Value = 99636.37  # L132: 99076.84 -> 99636.37
$ws.Cells.Item(132, 14).Value = -109756.37  # N132: -109196.84 -> -109756.37

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4170483.2  # H31: 3909956 -> 4170483.2
$ws.Cells.Item(31, 9).Value = 2466.1667  # I31: 2299.5715 -> 2466.1667
$ws.Cells.Item(31, 10).Value = 5212487.5  # J31: 5004100 -> 5212487.5
$ws.Cells.Item(31, 11).Value = 2466.1667  # K31: 2299.5715 -> 2466.1667
$ws.Cells.Item(31, 12).Value = 5212487.5  # L31: 5004100 -> 5212487.5
$ws.Cells.Item(31, 13).Value = -2171.1667  # M31: -2004.5715 -> -2171.1667
$ws.Cells.Item(31, 14).Value = -5213077.5  # N31: -5004690 -> -5213077.5
$ws.Cells.Item(34, 8).Value = 4170483.2  # H34: 3909956 -> 4170483.2
$ws.Cells.Item(34, 9).Value = 2466.1667  # I34: 2299.5715 -> 2466.1667
$ws.Cells.Item(34, 10).Value = 5212487.5  # J34: 5004100 -> 5212487.5
$ws.Cells.Item(34, 11).Value = 2466.1667  # K34: 2299.5715 -> 2466.1667
$ws.Cells.Item(34, 12).Value = 5212487.5  # L34: 5004100 -> 5212487.5
$ws.Cells.Item(34, 13).Value = -2264.1667  # M34: -2097.5715 -> -2264.1667
$ws.Cells.Item(34, 14).Value = -5212891.5  # N34: -5004504 -> -5212891.5
$ws.Cells.Item(58, 8).Value = 4516.8  # H58: 4024.1428 -> 4516.8
$ws.Cells.Item(58, 9).Value = 2792  # I58: 2792.25 -> 2792
$ws.Cells.Item(58, 11).Value = 2792  # K58: 2792.25 -> 2792
$ws.Cells.Item(58, 13).Value = -2589  # M58: -2589.25 -> -2589
$ws.Cells.Item(99, 8).Value = 4239.5  # H99: 4251.636 -> 4239.5
$ws.Cells.Item(99, 9).Value = 3989.2856  # I99: 3784.375 -> 3989.2856
$ws.Cells.Item(99, 10).Value = 4589.8  # J99: 5497.6665 -> 4589.8
$ws.Cells.Item(99, 11).Value = 3989.2856  # K99: 3784.375 -> 3989.2856
$ws.Cells.Item(99, 12).Value = 4589.8  # L99: 5497.6665 -> 4589.8
$ws.Cells.Item(99, 13).Value = -2491.2856  # M99: -2286.375 -> -2491.2856
$ws.Cells.Item(99, 14).Value = -7585.8  # N99: -8493.666499999999 -> -7585.8
$ws.Cells.Item(105, 8).Value = 1293.9697  # H105: 1248.1714 -> 1293.9697
$ws.Cells.Item(105, 9).Value = 1178.1562  # I105: 1137.8235 -> 1178.1562
$ws.Cells.Item(105, 11).Value = 1178.1562  # K105: 1137.8235 -> 1178.1562
$ws.Cells.Item(105, 13).Value = 568.8438000000001  # M105: 609.1765 -> 568.8438000000001
$ws.Cells.Item(126, 8).Value = 4239.5  # H126: 4251.636 -> 4239.5
$ws.Cells.Item(126, 9).Value = 3989.2856  # I126: 3784.375 -> 3989.2856
$ws.Cells.Item(126, 10).Value = 4589.8  # J126: 5497.6665 -> 4589.8
$ws.Cells.Item(126, 11).Value = 11967.8568  # K126: 11353.125 -> 11967.8568
$ws.Cells.Item(126, 12).Value = 13769.4  # L126: 16492.9995 -> 13769.4
$ws.Cells.Item(126, 13).Value = -9497.856800000001  # M126: -8883.125 -> -9497.856800000001
$ws.Cells.Item(126, 14).Value = -18709.4  # N126: -21432.9995 -> -18709.4
$ws.Cells.Item(132, 8).Value = 2667.875  # H132: 2318.5908 -> 2667.875
$ws.Cells.Item(132, 9).Value = 2044.1538  # I132: 1836.6842 -> 2044.1538
$ws.Cells.Item(132, 11).Value = 6132.4614  # K132: 5510.0526 -> 6132.4614
$ws.Cells.Item(132, 13).Value = -3602.4614  # M132: -2980.0526 -> -3602.4614
$ws.Cells.Item(134, 8).Value = 2459.5557  # H134: 2508.4285 -> 2459.5557
$ws.Cells.Item(134, 9).Value = 2319.5151  # I134: 2368.5938 -> 2319.5151
$ws.Cells.Item(134, 11).Value = 6958.5453  # K134: 7105.7814 -> 6958.5453
$ws.Cells.Item(134, 13).Value = -4423.5453  # M134: -4570.7814 -> -4423.5453
$ws.Cells.Item(136, 8).Value = 4516.8  # H136: 4024.1428 -> 4516.8
$ws.Cells.Item(136, 9).Value = 2792  # I136: 2792.25 -> 2792
$ws.Cells.Item(136, 11).Value = 8376  # K136: 8376.75 -> 8376
$ws.Cells.Item(136, 13).Value = -5826  # M136: -5826.75 -> -5826
$ws.Cells.Item(141, 8).Value = 140854.72  # H141: 157248.5 -> 140854.72
$ws.Cells.Item(141, 10).Value = 145997.17  # J141: 172998 -> 145997.17
$ws.Cells.Item(141, 12).Value = 145997.17  # L141: 172998 -> 145997.17
$ws.Cells.Item(141, 14).Value = -156357.17  # N141: -183358 -> -156357.17

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 357.3  # H34: 358.1 -> 357.3
$ws.Cells.Item(34, 9).Value = 321.75  # I34: 322.75 -> 321.75
$ws.Cells.Item(34, 11).Value = 965.25  # K34: 968.25 -> 965.25
$ws.Cells.Item(34, 13).Value = -881.25  # M34: -884.25 -> -881.25
$ws.Cells.Item(39, 8).Value = 6952.8335  # H39: 6143.6 -> 6952.8335
$ws.Cells.Item(39, 10).Value = 10116.875  # J39: 9822.833000000001 -> 10116.875
$ws.Cells.Item(39, 12).Value = 30350.625  # L39: 29468.499 -> 30350.625
$ws.Cells.Item(39, 14).Value = -30938.625  # N39: -30056.499 -> -30938.625
$ws.Cells.Item(40, 8).Value = 240.6  # H40: 306.8 -> 240.6
$ws.Cells.Item(40, 9).Value = 240.6  # I40: 306.8 -> 240.6
$ws.Cells.Item(40, 11).Value = 962.4  # K40: 1227.2 -> 962.4
$ws.Cells.Item(40, 13).Value = -893.4  # M40: -1158.2 -> -893.4
$ws.Cells.Item(55, 8).Value = 5518.278  # H55: 4583.1875 -> 5518.278
$ws.Cells.Item(55, 10).Value = 8560.1  # J55: 7450.375 -> 8560.1
$ws.Cells.Item(55, 12).Value = 25680.3  # L55: 22351.125 -> 25680.3
$ws.Cells.Item(55, 14).Value = -26034.3  # N55: -22705.125 -> -26034.3
$ws.Cells.Item(68, 8).Value = 3128962  # H68: 2781419.5 -> 3128962
$ws.Cells.Item(68, 9).Value = 1566.7142  # I68: 1495 -> 1566.7142
$ws.Cells.Item(68, 10).Value = 4004632.8  # J68: 3575683.8 -> 4004632.8
$ws.Cells.Item(68, 11).Value = 4700.142599999999  # K68: 4485 -> 4700.142599999999
$ws.Cells.Item(68, 12).Value = 12013898.4  # L68: 10727051.4 -> 12013898.4
$ws.Cells.Item(68, 13).Value = -3889.142599999999  # M68: -3674 -> -3889.142599999999
$ws.Cells.Item(68, 14).Value = -12015520.4  # N68: -10728673.4 -> -12015520.4
$ws.Cells.Item(69, 8).Value = 0  # H69: 1000 -> 0
$ws.Cells.Item(69, 9).Value = 0  # I69: 1000 -> 0
$ws.Cells.Item(69, 11).Value = 0  # K69: 3000 -> 0
$ws.Cells.Item(69, 13).ClearContents()  # M69: -2189 -> (removed)
$ws.Cells.Item(71, 8).Value = 3128962  # H71: 2781419.5 -> 3128962
$ws.Cells.Item(71, 9).Value = 1566.7142  # I71: 1495 -> 1566.7142
$ws.Cells.Item(71, 10).Value = 4004632.8  # J71: 3575683.8 -> 4004632.8
$ws.Cells.Item(71, 11).Value = 14100.4278  # K71: 13455 -> 14100.4278
$ws.Cells.Item(71, 12).Value = 36041695.2  # L71: 32181154.2 -> 36041695.2
$ws.Cells.Item(71, 13).Value = -10044.4278  # M71: -9399 -> -10044.4278
$ws.Cells.Item(71, 14).Value = -36049807.2  # N71: -32189266.2 -> -36049807.2
$ws.Cells.Item(72, 8).Value = 0  # H72: 1000 -> 0
$ws.Cells.Item(72, 9).Value = 0  # I72: 1000 -> 0
$ws.Cells.Item(72, 11).Value = 0  # K72: 9000 -> 0
$ws.Cells.Item(72, 13).ClearContents()  # M72: -4944 -> (removed)
$ws.Cells.Item(80, 8).Value = 2281  # H80: 1893.2 -> 2281
$ws.Cells.Item(80, 9).Value = 2281  # I80: 2153.6667 -> 2281
$ws.Cells.Item(80, 10).Value = 0  # J80: 1502.5 -> 0
$ws.Cells.Item(80, 11).Value = 6843  # K80: 6461.000100000001 -> 6843
$ws.Cells.Item(80, 12).Value = 0  # L80: 4507.5 -> 0
$ws.Cells.Item(80, 13).Value = -5907  # M80: -5525.000100000001 -> -5907
$ws.Cells.Item(80, 14).ClearContents()  # N80: -6379.5 -> (removed)
$ws.Cells.Item(83, 8).Value = 2281  # H83: 1893.2 -> 2281
$ws.Cells.Item(83, 9).Value = 2281  # I83: 2153.6667 -> 2281
$ws.Cells.Item(83, 10).Value = 0  # J83: 1502.5 -> 0
$ws.Cells.Item(83, 11).Value = 20529  # K83: 19383.0003 -> 20529
$ws.Cells.Item(83, 12).Value = 0  # L83: 13522.5 -> 0
$ws.Cells.Item(83, 13).Value = -15849  # M83: -14703.0003 -> -15849
$ws.Cells.Item(83, 14).ClearContents()  # N83: -22882.5 -> (removed)

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 8882.772000000001  # H102: 8192.916999999999 -> 8882.772000000001
$ws.Cells.Item(102, 9).Value = 4321.05  # I102: 3983.182 -> 4321.05
$ws.Cells.Item(102, 11).Value = 4321.05  # K102: 3983.182 -> 4321.05
$ws.Cells.Item(102, 13).Value = -2699.05  # M102: -2361.182 -> -2699.05
$ws.Cells.Item(126, 8).Value = 5313.364  # H126: 4860.923 -> 5313.364
$ws.Cells.Item(126, 9).Value = 4271.8887  # I126: 4079.2 -> 4271.8887
$ws.Cells.Item(126, 10).Value = 10000  # J126: 7466.6665 -> 10000
$ws.Cells.Item(126, 11).Value = 12815.6661  # K126: 12237.6 -> 12815.6661
$ws.Cells.Item(126, 12).Value = 30000  # L126: 22399.9995 -> 30000
$ws.Cells.Item(126, 13).Value = -10345.6661  # M126: -9767.599999999999 -> -10345.6661
$ws.Cells.Item(126, 14).Value = -34940  # N126: -27339.9995 -> -34940
$ws.Cells.Item(132, 8).Value = 2203.7778  # H132: 2202.6445 -> 2203.7778
$ws.Cells.Item(132, 9).Value = 2163.2666  # I132: 2161.5667 -> 2163.2666
$ws.Cells.Item(132, 11).Value = 6489.7998  # K132: 6484.7001 -> 6489.7998
$ws.Cells.Item(132, 13).Value = -3959.7998  # M132: -3954.7001 -> -3959.7998

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2574.5  # H7: 2799.4285 -> 2574.5
$ws.Cells.Item(7, 9).Value = 2574.5  # I7: 2599.3333 -> 2574.5
$ws.Cells.Item(7, 10).Value = 0  # J7: 4000 -> 0
$ws.Cells.Item(7, 11).Value = 2574.5  # K7: 2599.3333 -> 2574.5
$ws.Cells.Item(7, 12).Value = 0  # L7: 4000 -> 0
$ws.Cells.Item(7, 13).Value = -2462.5  # M7: -2487.3333 -> -2462.5
$ws.Cells.Item(7, 14).ClearContents()  # N7: -4224 -> (removed)
$ws.Cells.Item(40, 8).Value = 50079.54  # H40: 61296.617 -> 50079.54
$ws.Cells.Item(40, 9).Value = 58003.137  # I40: 78420.625 -> 58003.137
$ws.Cells.Item(40, 10).Value = 6499.75  # J40: 6499.8 -> 6499.75
$ws.Cells.Item(40, 11).Value = 58003.137  # K40: 78420.625 -> 58003.137
$ws.Cells.Item(40, 12).Value = 6499.75  # L40: 6499.8 -> 6499.75
$ws.Cells.Item(40, 13).Value = -57867.137  # M40: -78284.625 -> -57867.137
$ws.Cells.Item(40, 14).Value = -6771.75  # N40: -6771.8 -> -6771.75
$ws.Cells.Item(68, 8).Value = 2044.4  # H68: 2657.3333 -> 2044.4
$ws.Cells.Item(68, 9).Value = 2044.4  # I68: 2657.3333 -> 2044.4
$ws.Cells.Item(68, 11).Value = 2044.4  # K68: 2657.3333 -> 2044.4
$ws.Cells.Item(68, 13).Value = -1295.4  # M68: -1908.3333 -> -1295.4
$ws.Cells.Item(71, 8).Value = 2044.4  # H71: 2657.3333 -> 2044.4
$ws.Cells.Item(71, 9).Value = 2044.4  # I71: 2657.3333 -> 2044.4
$ws.Cells.Item(71, 11).Value = 10222  # K71: 13286.6665 -> 10222
$ws.Cells.Item(71, 13).Value = -6478  # M71: -9542.666499999999 -> -6478
$ws.Cells.Item(122, 8).Value = 2793.3809  # H122: 2982 -> 2793.3809
$ws.Cells.Item(122, 9).Value = 1267.5333  # I122: 1337.1666 -> 1267.5333
$ws.Cells.Item(122, 10).Value = 6608  # J122: 6929.6 -> 6608
$ws.Cells.Item(122, 11).Value = 3802.5999  # K122: 4011.4998 -> 3802.5999
$ws.Cells.Item(122, 12).Value = 19824  # L122: 20788.8 -> 19824
$ws.Cells.Item(122, 13).Value = -1352.5999  # M122: -1561.4998 -> -1352.5999
$ws.Cells.Item(122, 14).Value = -24724  # N122: -25688.8 -> -24724
$ws.Cells.Item(126, 8).Value = 2574.5  # H126: 2799.4285 -> 2574.5
$ws.Cells.Item(126, 9).Value = 2574.5  # I126: 2599.3333 -> 2574.5
$ws.Cells.Item(126, 10).Value = 0  # J126: 4000 -> 0
$ws.Cells.Item(126, 11).Value = 7723.5  # K126: 7797.999899999999 -> 7723.5
$ws.Cells.Item(126, 12).Value = 0  # L126: 12000 -> 0
$ws.Cells.Item(126, 13).Value = -5253.5  # M126: -5327.999899999999 -> -5253.5
$ws.Cells.Item(126, 14).ClearContents()  # N126: -16940 -> (removed)
$ws.Cells.Item(132, 8).Value = 7463.722  # H132: 6006 -> 7463.722
$ws.Cells.Item(132, 9).Value = 7804.364  # I132: 5647.278 -> 7804.364
$ws.Cells.Item(132, 11).Value = 23413.092  # K132: 16941.834 -> 23413.092
$ws.Cells.Item(132, 13).Value = -20883.092  # M132: -14411.834 -> -20883.092

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 9999.5  # H62: 3999 -> 9999.5
$ws.Cells.Item(62, 9).Value = 0  # I62: 3999 -> 0
$ws.Cells.Item(62, 10).Value = 9999.5  # J62: 0 -> 9999.5
$ws.Cells.Item(62, 11).Value = 0  # K62: 3999 -> 0
$ws.Cells.Item(62, 12).Value = 9999.5  # L62: 0 -> 9999.5
$ws.Cells.Item(62, 13).ClearContents()  # M62: -3375 -> (removed)
$ws.Cells.Item(62, 14).Value = -11247.5  # N62: None -> -11247.5
$ws.Cells.Item(65, 8).Value = 9999.5  # H65: 3999 -> 9999.5
$ws.Cells.Item(65, 9).Value = 0  # I65: 3999 -> 0
$ws.Cells.Item(65, 10).Value = 9999.5  # J65: 0 -> 9999.5
$ws.Cells.Item(65, 11).Value = 0  # K65: 19995 -> 0
$ws.Cells.Item(65, 12).Value = 49997.5  # L65: 0 -> 49997.5
$ws.Cells.Item(65, 13).ClearContents()  # M65: -16875 -> (removed)
$ws.Cells.Item(65, 14).Value = -56237.5  # N65: None -> -56237.5
$ws.Cells.Item(81, 8).Value = 5601.72  # H81: 6252.273 -> 5601.72
$ws.Cells.Item(81, 9).Value = 6235.1  # I81: 6872.222 -> 6235.1
$ws.Cells.Item(81, 10).Value = 5179.467  # J81: 5823.077 -> 5179.467
$ws.Cells.Item(81, 11).Value = 12470.2  # K81: 13744.444 -> 12470.2
$ws.Cells.Item(81, 12).Value = 10358.934  # L81: 11646.154 -> 10358.934
$ws.Cells.Item(81, 13).Value = -11409.2  # M81: -12683.444 -> -11409.2
$ws.Cells.Item(81, 14).Value = -12480.934  # N81: -13768.154 -> -12480.934
$ws.Cells.Item(84, 8).Value = 5601.72  # H84: 6252.273 -> 5601.72
$ws.Cells.Item(84, 9).Value = 6235.1  # I84: 6872.222 -> 6235.1
$ws.Cells.Item(84, 10).Value = 5179.467  # J84: 5823.077 -> 5179.467
$ws.Cells.Item(84, 11).Value = 62351  # K84: 68722.22 -> 62351
$ws.Cells.Item(84, 12).Value = 51794.67  # L84: 58230.77 -> 51794.67
$ws.Cells.Item(84, 13).Value = -57047  # M84: -63418.22 -> -57047
$ws.Cells.Item(84, 14).Value = -62402.67  # N84: -68838.77 -> -62402.67
$ws.Cells.Item(126, 8).Value = 29747.666  # H126: 29945.223 -> 29747.666
$ws.Cells.Item(126, 9).Value = 32841.125  # I126: 42417.168 -> 32841.125
$ws.Cells.Item(126, 10).Value = 5000  # J126: 5001.3335 -> 5000
$ws.Cells.Item(126, 11).Value = 98523.375  # K126: 127251.504 -> 98523.375
$ws.Cells.Item(126, 12).Value = 15000  # L126: 15004.0005 -> 15000
$ws.Cells.Item(126, 13).Value = -96053.375  # M126: -124781.504 -> -96053.375
$ws.Cells.Item(126, 14).Value = -19940  # N126: -19944.0005 -> -19940
$ws.Cells.Item(132, 8).Value = 3475.8333  # H132: 3613.3823 -> 3475.8333
$ws.Cells.Item(132, 9).Value = 3213.4688  # I132: 3351.8667 -> 3213.4688
$ws.Cells.Item(132, 11).Value = 9640.4064  # K132: 10055.6001 -> 9640.4064
$ws.Cells.Item(132, 13).Value = -7110.4064  # M132: -7525.6001 -> -7110.4064
